# Automatische test-sync: 2025-06-20 10:30:50
# Adds the new "Afmelding nieuwsbrief" unsubscribe log entry to the Logs
# sheet and refreshes the Dashboard summary sheet to match.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append the new mail log entry as row 7 -------------------
$logs.Range("A7").Value = "Afmelding nieuwsbrief"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D7").Value = "Afmelding / Nieuwsbrief"
$logs.Range("F7").Value = "2025-06-20 10:30:11"
$logs.Range("G7").Value = "Nee"

# --- Dashboard sheet: refresh summary counts -------------------------------
# New entry bumped "Afmelding / Nieuwsbrief" to 2, tying with
# "Samenwerking / Partnerverzoek" (also 2) -- so the two rows swap order.
$dashboard.Range("A2").Value = "Afmelding / Nieuwsbrief"
$dashboard.Range("B2").Value = 2

$dashboard.Range("A3").Value = "Samenwerking / Partnerverzoek"
$dashboard.Range("B3").Value = 2

# --- Logs sheet: extend conditional formatting ranges to cover row 7 ------
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D7"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G7"))
